$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was collected for this product/market and needs
# to be inserted as the new first row of the Naranja (Valencia) data block,
# pushing every existing record down by one row (row 670 -> 671, ...,
# 766 -> 767), exactly like typing a new row at the top of a running log.
$ws.Rows.Item(670).Insert()

# Populate the newly inserted row 670 with the new observation.
$ws.Cells.Item(670, 1).Value = 5
$ws.Cells.Item(670, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(670, 3).Value = "Maule"
$ws.Cells.Item(670, 4).Value = 44951
$ws.Cells.Item(670, 5).Value = 7
$ws.Cells.Item(670, 6).Value = "Fruta"
$ws.Cells.Item(670, 7).Value = 100102
$ws.Cells.Item(670, 8).Value = "Cítricos"
$ws.Cells.Item(670, 9).Value = 100102005
$ws.Cells.Item(670, 10).Value = "Naranja"
$ws.Cells.Item(670, 11).Value = "Valencia"
$ws.Cells.Item(670, 12).Value = "Primera"
$ws.Cells.Item(670, 13).Value = 400
$ws.Cells.Item(670, 14).Value = 11000
$ws.Cells.Item(670, 15).Value = 11000
$ws.Cells.Item(670, 16).Value = 11000
$ws.Cells.Item(670, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(670, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(670, 19).Value = 733
$ws.Cells.Item(670, 20).Value = 15
